$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend Table1 to include the new Country column (B3:E7)
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("B3:E7"))

# Fill in the new column's header and data
$ws.Range("E3").Value = "Country"
$ws.Range("E4").Value = "India"
$ws.Range("E5").Value = "India"
